$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.516.53"
$ws.Range("E2").Value = "  +0.46%  "

$ws.Range("D3").Value = "3.111.12"
$ws.Range("E3").Value = "  +0.50%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "525.62"
$ws.Range("E5").Value = "  +0.60%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "136.86"
$ws.Range("E6").Value = "  -2.72%  "

$ws.Range("E7").Value = "  +0.05%  "

$ws.Range("D8").Value = "3.106.66"
$ws.Range("E8").Value = "  +0.37%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.448"
$ws.Range("E9").Value = "  +2.52%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.25"
$ws.Range("E10").Value = "  +0.95%  "

$ws.Range("E11").Value = "  -0.22%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.397"
$ws.Range("E12").Value = "  +3.16%  "

$ws.Range("D13").Value = "3.642.38"
$ws.Range("E13").Value = "  +0.38%  "

$ws.Range("E14").Value = "  +3.01%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "25.34"
$ws.Range("E15").Value = "  -2.53%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000164"
$ws.Range("E16").Value = "  +0.59%  "

$ws.Range("D17").Value = "57.589.54"
$ws.Range("E17").Value = "  +0.46%  "

$ws.Range("D18").Value = "3.108.12"
$ws.Range("E18").Value = "  +0.30%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.95"
$ws.Range("E19").Value = "  -2.37%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.47"
$ws.Range("E20").Value = "  -2.37%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.89"
$ws.Range("E21").Value = "  -1.73%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "348.34"
$ws.Range("E22").Value = "  +3.32%  "

$ws.Range("E23").Value = "  -0.12%  "

$ws.Range("E24").Value = "  +0.24%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "68.16"
$ws.Range("E25").Value = "  +2.33%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.503"
$ws.Range("E26").Value = "  -1.73%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.168"
$ws.Range("E27").Value = "  +0.03%  "

$ws.Range("E28").Value = "  +0.05%  "

$ws.Range("D29").Value = "0.0₃0909"
$ws.Range("E29").Value = "  +0.16%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.43"
$ws.Range("E30").Value = "  +3.87%  "

$ws.Range("E31").Value = "  +0.06%  "

$ws.Range("E32").Value = "  +0.62%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.00"
$ws.Range("E33").Value = "  -7.30%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "20.98"
$ws.Range("E34").Value = "  +0.31%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.98"
$ws.Range("E35").Value = "  +7.91%  "

$ws.Range("E36").Value = "  -1.86%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "158.56"
$ws.Range("E37").Value = "  +0.95%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.11"
$ws.Range("E38").Value = "  +0.35%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "25.91"
$ws.Range("E39").Value = "  -4.16%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.24"
$ws.Range("E40").Value = "  -2.87%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.22"
$ws.Range("E41").Value = "  +7.15%  "

$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.62"
$ws.Range("E42").Value = "  +6.89%  "

$ws.Range("B43").Value = "Hedera"
$ws.Range("C43").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0664"
$ws.Range("E43").Value = "  +1.23%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.699"
$ws.Range("E44").Value = "  +2.19%  "

$ws.Range("D45").Value = "3.144.56"
$ws.Range("E45").Value = "  +0.18%  "

$ws.Range("B46").Value = "Maker"
$ws.Range("C46").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D46").Value = "2.354.27"
$ws.Range("E46").Value = "  +2.14%  "

$ws.Range("B47").Value = "OKB"
$ws.Range("C47").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "36.49"
$ws.Range("E47").Value = "  -0.28%  "

$ws.Range("E48").Value = "  -0.06%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0268"
$ws.Range("E49").Value = "  +3.52%  "

$ws.Range("B50").Value = "ONDO"
$ws.Range("C50").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.957"
$ws.Range("E50").Value = "  -1.09%  "

$ws.Range("B51").Value = "Cosmos"
$ws.Range("C51").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.04"
$ws.Range("E51").Value = "  +0.61%  "
